$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 240; this shifts rows 240-279 down to 241-280
# (Excel copies formatting from the row above, matching the style shift seen
# in the diff where D240's date style moves down with the data.)
$ws.Rows.Item(240).Insert()

# Populate the newly inserted row 240 with the new record's data.
$ws.Range("A240").Value = 6
$ws.Range("B240").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C240").Value = "Metropolitana"
$ws.Range("D240").Value = 44951
$ws.Range("E240").Value = 13
$ws.Range("F240").Value = 100112001
$ws.Range("G240").Value = "Berenjena"
$ws.Range("H240").Value = "Sin especificar"
$ws.Range("I240").Value = "Primera"
$ws.Range("J240").Value = 180
$ws.Range("K240").Value = 18000
$ws.Range("L240").Value = 20000
$ws.Range("M240").Value = 18889
$ws.Range("N240").Value = "$/caja 50 unidades"
$ws.Range("O240").Value = "Región Metropolitana"
$ws.Range("P240").Value = 378
$ws.Range("Q240").Value = 50
$ws.Range("R240").Value = "Hortaliza"
